$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> FAPs (self loop), refreshed TPM-derived values
$ws.Range("D2").Value = "FAPs"
$ws.Range("H2").Value = 0.779989
$ws.Range("M2").Value = 0.1534146666666667
$ws.Range("N2").Value = 0.460244
$ws.Range("O2").Value = 0.8701728646218362
$ws.Range("P2").Value = 0.8701728646218362
$ws.Range("Q2").Value = 0.03988725081288889
$ws.Range("R2").Value = 0.358985257316
$ws.Range("S2").Value = 0.8701728646218362
$ws.Range("T2").Value = 0.8701728646218362

# Row 3: FAPs -> Resolving-Mac, refreshed TPM-derived values
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("H3").Value = 0.779989
$ws.Range("M3").Value = 0.022889
$ws.Range("N3").Value = 0.06866700000000001
$ws.Range("O3").Value = 0.1298271353781638
$ws.Range("P3").Value = 0.1298271353781638
$ws.Range("Q3").Value = 0.005951056073666667
$ws.Range("R3").Value = 0.05355950466300001
$ws.Range("S3").Value = 0.1298271353781638
$ws.Range("T3").Value = 0.1298271353781638

# Remove now-obsolete target-cluster rows (Inflammatory-Mac, MuSCs, and the old
# trailing Resolving-Mac row) so only two data rows remain under the header.
$ws.Rows("4:6").Delete()
